# Applies the "add locks for lt cancel and update report" edit:
#  - Rename header labels for the three "Lifetime Cancels" columns to
#    "New Lifetime Cancels" variants.
#  - Update a handful of recalculated metric values across several rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Header label updates (shared strings) ---
$ws.Range("D1").Value = "New Lifetime Cancels"
$ws.Range("N1").Value = "New Lifetime Cancels DVH"
$ws.Range("S1").Value = "New Lifetime Cancels Copay"

# --- Updated metric values ---
$ws.Range("C3").Value  = 975
$ws.Range("C7").Value  = 953
$ws.Range("C10").Value = 492
$ws.Range("J10").Value = 28
$ws.Range("C11").Value = 455
$ws.Range("J11").Value = 1
$ws.Range("C13").Value = 111
$ws.Range("J13").Value = 0
$ws.Range("C14").Value = 154
$ws.Range("C15").Value = 736
$ws.Range("C16").Value = 884
$ws.Range("J16").Value = 394
$ws.Range("C18").Value = 747
$ws.Range("J18").Value = 144
$ws.Range("C19").Value = 623
$ws.Range("J19").Value = 55
$ws.Range("C21").Value = 1143
$ws.Range("J21").Value = 517
$ws.Range("M22").Value = 59
